$d = $word.ActiveDocument

# 1. Title: Appendix 12 -> Bylaag 12
$d.Content.Find.Execute(
    "Appendix 12: SWIFT Child Safety Module Interview: Participants", $true, $false, $false, $false, $false,
    $true, 1, $false, "Bylaag 12: SWIFT-kinderveiligheidsmodule Onderhoud: Deelnemers", 2)

# 2. "Wat het jy nie ..." -> "Waarvan het jy nie ..."
$d.Content.Find.Execute(
    "Wat het jy nie van hierdie module gehou nie? Hoe kan ons ", $true, $false, $false, $false, $false,
    $true, 1, $false, "Waarvan het jy nie van hierdie module gehou nie? Hoe kan ons ", 2)

# 3. "What stands out for you about the " -> "Wat staan vir jou uit oor die "
$d.Content.Find.Execute(
    "What stands out for you about the ", $true, $false, $false, $false, $false,
    $true, 1, $false, "Wat staan vir jou uit oor die ", 2)

# 4. "online safety" -> "aanlyn-veiligheid"
$d.Content.Find.Execute(
    "online safety", $true, $false, $false, $false, $false,
    $true, 1, $false, "aanlyn-veiligheid", 2)

# 5. The lone space before "lesson" becomes a hyphen (keeps run boundaries intact)
$r = $d.Content
$r.Find.Execute("’ lesson", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$spaceRange = $d.Range($r.Start + 1, $r.Start + 2)
$spaceRange.Text = "-"

# 6. "lesson" -> "les"
$d.Content.Find.Execute(
    "lesson", $true, $false, $false, $false, $false,
    $true, 1, $false, "les", 2)

# 7. "Would you mind sharing with me what that was?" -> "Sal jy omgee om met my te deel wat dit was?"
$d.Content.Find.Execute(
    "Would you mind sharing with me what that was?", $true, $false, $false, $false, $false,
    $true, 1, $false, "Sal jy omgee om met my te deel wat dit was?", 2)
